$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DeviceList")

# Delete columns F and H (device columns to remove)
$ws.Range("H1").EntireColumn.Delete() | Out-Null
$ws.Range("F1").EntireColumn.Delete() | Out-Null

# Fix Individual_ID row (row 9) back to sequential 1..8
$ws.Range("F9").Value2 = 5
$ws.Range("G9").Value2 = 6
$ws.Range("H9").Value2 = 7
$ws.Range("I9").Value2 = 8

for ($r = 1; $r -le 10; $r++) {
  $line = "Row " + $r + ": "
  for ($c = 1; $c -le 9; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $v = $cell.Value2
    $line += "[" + $v + "] "
  }
  Write-Host $line
}
